$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 69 (shifts existing rows 69-178 down to 70-179)
$ws.Rows(69).Insert()

# Populate the newly inserted row 69 with the new weekly record
$ws.Range("A69").Value = 4
$ws.Range("B69").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C69").Value = "Los Lagos"
$ws.Range("D69").Value = 45219
$ws.Range("E69").Value = 10
$ws.Range("F69").Value = 100112031
$ws.Range("G69").Value = "Poroto verde"
$ws.Range("H69").Value = "Magnum"
$ws.Range("I69").Value = "Primera"
$ws.Range("J69").Value = 40
$ws.Range("K69").Value = 37000
$ws.Range("L69").Value = 37000
$ws.Range("M69").Value = 37000
$ws.Range("N69").Value = "`$/malla 25 kilos"
$ws.Range("O69").Value = "Perú"
$ws.Range("P69").Value = 1480
$ws.Range("Q69").Value = 25
$ws.Range("R69").Value = "Hortaliza"

# Match the original date-cell number format used throughout column D
$ws.Range("D69").NumberFormat = $ws.Range("D70").NumberFormat
